# Apply updated TPM values to rows 2-7, then remove now-obsolete rows 8-10 (ECs target-cluster rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("G2").Value2 = 12.056684
$ws.Range("H2").Value2 = 36.170052
$ws.Range("I2").Value2 = 0.06307822458376462
$ws.Range("J2").Value2 = 0.06307822458376462
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.2109236666666666
$ws.Range("N2").Value2 = 0.632771
$ws.Range("O2").Value2 = 0.8951984155054113
$ws.Range("P2").Value2 = 0.8951984155054113
$ws.Range("Q2").Value2 = 2.543039997121333
$ws.Range("R2").Value2 = 22.887359974092
$ws.Range("S2").Value2 = 0.05646752670028057
$ws.Range("T2").Value2 = 0.05646752670028057

# Row 3
$ws.Range("D3").Value2 = "MuSCs"
$ws.Range("G3").Value2 = 12.056684
$ws.Range("H3").Value2 = 36.170052
$ws.Range("I3").Value2 = 0.06307822458376462
$ws.Range("J3").Value2 = 0.06307822458376462
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.024693
$ws.Range("N3").Value2 = 0.07407900000000001
$ws.Range("O3").Value2 = 0.1048015844945887
$ws.Range("P3").Value2 = 0.1048015844945887
$ws.Range("Q3").Value2 = 0.297715698012
$ws.Range("R3").Value2 = 2.679441282108
$ws.Range("S3").Value2 = 0.006610697883484049
$ws.Range("T3").Value2 = 0.006610697883484049

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("G4").Value2 = 25.140634
$ws.Range("H4").Value2 = 75.421902
$ws.Range("I4").Value2 = 0.1315309049843414
$ws.Range("J4").Value2 = 0.1315309049843414
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.2109236666666666
$ws.Range("N4").Value2 = 0.632771
$ws.Range("O4").Value2 = 0.8951984155054113
$ws.Range("P4").Value2 = 0.8951984155054113
$ws.Range("Q4").Value2 = 5.302754705604666
$ws.Range("R4").Value2 = 47.724792350442
$ws.Range("S4").Value2 = 0.1177462577319752
$ws.Range("T4").Value2 = 0.1177462577319752

# Row 5
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("I5").Value2 = 0.1315309049843414
$ws.Range("J5").Value2 = 0.1315309049843414
$ws.Range("M5").Value2 = 0.024693
$ws.Range("N5").Value2 = 0.07407900000000001
$ws.Range("O5").Value2 = 0.1048015844945887
$ws.Range("P5").Value2 = 0.1048015844945887
$ws.Range("Q5").Value2 = 0.6207976753620001
$ws.Range("R5").Value2 = 5.587179078258001
$ws.Range("S5").Value2 = 0.01378464725236617
$ws.Range("T5").Value2 = 0.01378464725236617

# Row 6
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("G6").Value2 = 153.9412893333333
$ws.Range("H6").Value2 = 461.8238680000001
$ws.Range("I6").Value2 = 0.8053908704318941
$ws.Range("J6").Value2 = 0.8053908704318941
$ws.Range("O6").Value2 = 0.8951984155054113
$ws.Range("P6").Value2 = 0.8951984155054113
$ws.Range("Q6").Value2 = 32.46986119758089
$ws.Range("R6").Value2 = 292.228750778228
$ws.Range("S6").Value2 = 0.7209846310731556
$ws.Range("T6").Value2 = 0.7209846310731556

# Row 7
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("G7").Value2 = 153.9412893333333
$ws.Range("H7").Value2 = 461.8238680000001
$ws.Range("I7").Value2 = 0.8053908704318941
$ws.Range("J7").Value2 = 0.8053908704318941
$ws.Range("N7").Value2 = 0.07407900000000001
$ws.Range("O7").Value2 = 0.1048015844945887
$ws.Range("P7").Value2 = 0.1048015844945887
$ws.Range("Q7").Value2 = 3.801272257508001
$ws.Range("R7").Value2 = 34.21145031757201
$ws.Range("S7").Value2 = 0.08440623935873848
$ws.Range("T7").Value2 = 0.08440623935873846

# Remove rows 8-10 (ECs as target cluster no longer present in filtered results)
$ws.Rows("8:10").Delete()
